{"js": "const QUESTIONS = [\"(1) \u738b\u8001\u5e2b\u53bb\u8cb7\u5152\u7ae5\u5c0f\u63d0\u7434\uff0c\u82e5\u8cb77\u628a\uff0c\u5247\u6240\u5e36\u7684\u9322\u5dee110\u5143\uff1b\u82e5\u8cb75\u628a\uff0c\u5247\u6240\u5e36\u7684\u9322\u9084\u5dee30\u5143\u3002\u554f\uff1a\u5152\u7ae5\u5c0f\u63d0\u7434\u4e00\u628a\u301040\u3011\u5143\uff0c\u738b\u8001\u5e2b\u5e36\u4e86\u3010170\u3011\u5143\u3002\", \"(2) \u4e00\u500b\u5546\u8ca9\u4f30\u8a08\uff0c\u82e51\u516c\u65a4\u860b\u679c\u8ce324\u5143\uff0c\u6703\u8ce040\u5143\uff1b\u82e51\u516c\u65a4\u860b\u679c\u8ce330\u5143\uff0c\u53ef\u4ee5\u8cfa80\u5143\u3002\u554f\uff1a\u82e5\u4ee5\u4e0d\u8ce0\u4e0d\u8cfa\u7684\u50f9\u683c\u8ce3\u51fa\uff0c\u6bcf\u516c\u65a4\u860b\u679c\u61c9\u8ce3\u301026\u3011\u5143\", \"(3) \u8c93\u5abd\u5abd\u7d66\u5c0f\u8c93\u5206\u9b5a\uff0c\u6bcf\u96bb\u5c0f\u8c93\u520610\u689d\u9b5a\uff0c\u5c31\u591a\u51fa8\u689d\u9b5a\uff0c\u6bcf\u96bb\u5c0f\u8c93\u520611\u689d\u9b5a\u5247\u6b63\u597d\u5206\u5b8c\u3002\u554f\uff1a\u4e00\u5171\u6709\u30108\u3011\u96bb\u5c0f\u8c93\u3001\u8c93\u5abd\u5abd\u4e00\u5171\u6709\u301088\u3011\u689d\u9b5a\", \"(4) \u5b78\u751f\u5011\u53c3\u52a0\u690d\u6a39\u6d3b\u52d5\uff0c\u5982\u679c\u6bcf\u4eba\u683d5\u68f5\u6a39\uff0c\u9084\u526912\u68f5\u6a39\uff1b\u5982\u679c\u6bcf\u4eba\u683d7\u68f5\uff0c\u5c31\u7f3a4\u68f5\u3002\u554f\uff1a\u5b78\u751f\u6709\u30108\u3011\u4eba\u3001\u4e00\u5171\u8981\u683d\u301052\u3011\u68f5\u6a39\", \"(5) 401\u73ed\u540c\u5b78\u690d\u6a39\uff0c\u6bcf\u4eba\u690d1\u68f5\u9084\u526920\u68f5\uff0c\u6bcf\u4eba\u690d2\u68f5\u5dee30\u68f5\u3002\u554f\uff1a\u6709\u301050\u3011\u500b\u540c\u5b78\u3001\u301070\u3011\u68f5\u6a39\u82d7\", \"(6) \u5abd\u5abd\u8cb7\u56de\u4e00\u7b50\u860b\u679c\uff0c\u6309\u8a08\u5283\u5403\u7684\u5929\u6578\u7b97\u4e86\u4e00\u4e0b\uff0c\u5982\u679c\u6bcf\u5929\u54034\u500b\uff0c\u8981\u591a\u51fa48\u500b\u860b\u679c\uff1b\u5982\u679c\u6bcf\u5929\u54036\u500b\uff0c\u5247\u53c8\u5c118\u500b\u860b\u679c\u3002\u554f\uff1a\u5abd\u5abd\u8cb7\u56de\u7684\u860b\u679c\u6709\u3010160\u3011\u500b\u3001\u8a08\u756b\u5403\u301028\u3011\u5929 \", \"(7) \u5c0f\u8ecd\u5c07\u81ea\u5df1\u6536\u85cf\u7684\u4e00\u4e9b\u7167\u7247\u9001\u7d66\u5e7c\u5152\u5712\u5927\u73ed\u7684\u5c0f\u670b\u53cb\u5011\u3002\u5982\u679c\u6bcf\u4eba\u52069\u5f35\uff0c\u9084\u591a12\u5f35\uff0c\u5982\u679c\u6bcf\u4eba\u520610\u5f35\u5247\u6b63\u597d\u5206\u5b8c\u3002\u554f\uff1a\u6709\u301012\u3011\u500b\u5c0f\u670b\u53cb\uff0c\u7167\u7247\u6709\u3010120\u3011\u5f35\", \"(8) \u4e00\u5305\u7cd6\u5206\u7d66\u5e7e\u500b\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u6bcf\u4eba\u52063\u584a\uff0c\u5247\u99183\u584a\uff1b\u5982\u679c\u6bcf\u4eba\u52065\u584a\uff0c\u5247\u5c111\u584a\u3002\u554f\uff1a\u5c0f\u670b\u53cb\u6709\u30102\u3011\u4eba\u3001\u7cd6\u6709\u30109\u3011\u584a\", \"(9) \u7528\u7b50\u88dd\u897f\u74dc\uff0c\u5982\u679c\u6bcf\u7b50\u88dd5\u500b\uff0c\u5247\u5c1115\u500b\u897f\u74dc\uff1b\u5982\u679c\u6bcf\u7b50\u88dd3\u500b\uff0c\u5247\u591a29\u500b\u897f\u74dc\u3002\u554f\uff1a\u5171\u6709\u7b50\u301022\u3011\u500b\u3001\u897f\u74dc\u301095\u3011\u500b\", \"(10) \u6709\u4e00\u6279\u7df4\u7fd2\u672c\u767c\u7d66\u5b78\u751f\uff0c\u5982\u679c\u6bcf\u4eba5\u672c\uff0c\u5247\u591a70\u672c\uff0c\u5982\u679c\u6bcf\u4eba7\u672c\uff0c\u5247\u591a10\u672c\u3002\u554f\uff1a\u9019\u500b\u73ed\u6709\u5b78\u751f\u301030\u3011\u4eba\uff0c\u6709\u7df4\u7fd2\u672c\u3010220\u3011\u672c\u3002\", \"(11) \u8001\u7334\u5b50\u7d66\u5c0f\u7334\u5b50\u5206\u6843\uff0c\u6bcf\u96bb\u5c0f\u7334\u520610\u500b\u6843\uff0c\u5c31\u591a\u51fa9\u500b\u6843\uff0c\u6bcf\u96bb\u5c0f\u7334\u520611\u500b\u6843\u5247\u591a\u51fa2\u500b\u6843\u3002\u554f\uff1a\u4e00\u5171\u6709\u3010 7 \u3011\u96bb\u5c0f\u7334\u5b50\uff0c\u8001\u7334\u5b50\u4e00\u5171\u6709\u3010 79 \u3011\u500b\u6843\u5b50\", \"(12) \u5b78\u6821\u5c07\u4e00\u6279\u925b\u7b46\u734e\u7d66\u5b78\u751f\uff0c\u6bcf\u4eba9\u652f\u7f3a15\u652f\uff1b\u6bcf\u4eba7\u652f\u7f3a7\u652f\u3002\u554f\uff1a\u5b78\u751f\u6709\u30104\u3011\u4eba\uff0c\u925b\u7b46\u6709\u301021\u3011\u652f\", \"(13) \u5c0f\u670b\u53cb\u5206\u7cd6\u679c\uff0c\u6bcf\u4eba3\u7c92\uff0c\u991830\u7c92\uff1b\u6bcf\u4eba5\u7c92\uff0c\u5c114\u7c92\u3002\u554f\uff1a\u6709\u301017\u3011\u500b\u5c0f\u670b\u53cb\u3001\u301081\u3011\u7c92\u7cd6\", \"(14) \u4e00\u500b\u6c7d\u8eca\u968a\u904b\u8f38\u4e00\u6279\u8ca8\u7269\uff0c\u5982\u679c\u6bcf\u8f1b\u6c7d\u8eca\u904b3500\u516c\u65a4\uff0c\u90a3\u9ebc\u8ca8\u7269\u9084\u5269\u4e0b5000\u516c\u65a4\uff1b\u5982\u679c\u6bcf\u8f1b\u6c7d\u8eca\u904b4000\u516c\u65a4\uff0c\u90a3\u9ebc\u8ca8\u7269\u9084\u5269\u4e0b500\u516c\u65a4\u3002\u554f\uff1a\u9019\u500b\u6c7d\u8eca\u968a\u6709\u30109\u3011\u8f1b\u6c7d\u8eca\u3001\u8981\u904b\u7684\u8ca8\u7269\u6709\u301036500\u3011\u516c\u65a4\", \"(15) \u5b78\u6821\u5206\u914d\u82e5\u5e72\u4eba\u64e6\u73bb\u7483\uff0c\u5176\u4e2d\u5169\u4eba\u5404\u64e64\u584a\uff0c\u5176\u9918\u5404\u64e65\u584a\uff0c\u5247\u991812\u584a\uff1b\u82e5\u6bcf\u4eba\u64e66\u584a\uff0c\u5247\u6b63\u597d\u64e6\u5b8c\u3002\u554f\uff1a\u64e6\u73bb\u7483\u7684\u6709\u301010\u3011\u4eba\u3001\u73bb\u7483\u6709\u301060\u3011\u584a\", \"(16) \u5c07\u4e00\u5806\u6843\u5b50\uff0c\u5e73\u5747\u5206\u7d66\u7334\u5b50\uff0c\u6bcf\u96bb\u7334\u5b50\u520610\u500b\uff0c\u6709\u5169\u96bb\u7334\u5b50\u6c92\u5206\u5230\uff0c\u7b2c\u4e8c\u6b21\u91cd\u5206\uff0c\u6bcf\u96bb\u7334\u5b508\u500b\u6843\u5b50\uff0c\u525b\u5de7\u5206\u5b8c\u3002\u554f\uff1a\u4e00\u5806\u6843\u5b50\u6709\u301080\u3011\u500b\u3001\u7334\u5b50\u6709\u301010\u3011\u96bb\", \"(17) \u8fb2\u6c11\u92e4\u8349\uff0c\u5176\u4e2d5\u4eba\u5404\u92e44\u755d\uff0c\u9918\u4e0b\u7684\u5404\u92e43\u755d\uff0c\u9019\u6a23\u5206\u914d\u6700\u5f8c\u9918\u4e0b26\u755d\uff1b\u5982\u679c\u5176\u4e2d3\u4eba\u6bcf\u4eba\u5404\u92e43\u755d\uff0c\u9918\u4e0b\u7684\u4eba\u5404\u92e45\u755d\uff0c\u6700\u5f8c\u9918\u4e0b3\u755d\u3002\u554f\uff1a\u8349\u5730\u9762\u7a4d\u301082\u3011\u755d\u3001\u92e4\u8349\u4eba\u6578\u301017\u3011\u4eba\", \"(18) \u8001\u5e2b\u628a\u4e00\u7c43\u860b\u679c\u5206\u7d66\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u6e1b\u5c11\u4e00\u540d\u540c\u5b78\uff0c\u6bcf\u500b\u540c\u5b78\u6b63\u597d\u5206\u5f975\u500b\uff1b\u5982\u679c\u589e\u52a0\u4e00\u500b\u540c\u5b78\uff0c\u6b63\u597d\u6bcf\u4eba\u5206\u5f974\u500b\uff0c\u6c42\u9019\u7c43\u860b\u679c\u4e00\u5171\u6709\u3010 40 \u3011\u500b\", \"(19) \u5321\u660e\u548c\u674e\u6b23\u7d04\u5b9a\u5728\u76f8\u540c\u7684\u6642\u9593\u5167\u505a\u5b8c\u540c\u6a23\u7684\u4e00\u672c\u6578\u5b78\u7df4\u7fd2\u984c\u3002\u5321\u660e\u8a08\u756b\u982d\u5169\u5468\u6bcf\u9031\u505a30\u9053\uff0c\u4ee5\u5f8c\u6bcf\u9031\u505a25\u9053\uff1b\u674e\u6b23\u8a08\u756b\u982d\u5169\u5468\u6bcf\u9031\u505a35\u9053\uff0c\u4ee5\u5f8c\u6bcf\u9031\u505a30\u9053\u3002\u7d50\u679c\u674e\u6b23\u63d0\u524d\u5169\u5468\u505a\u5b8c\u3002\u554f\uff1a\u4ed6\u5011\u6e96\u5099\u301012\u3011\u5468\u505a\u5b8c\uff0c\u9019\u672c\u7df4\u7fd2\u984c\u5171\u6709\u3010310\u3011\u984c\", \"(20) \u7ae5\u8ecd\u968a\u54e1\u64fa\u82b1\u76c6\u5e03\u7f6e\u6821\u5712\u3002\u5982\u679c\u6bcf\u4eba\u64fa5\u76c6\u82b1\uff0c\u9084\u67093\u76c6\u6c92\u4eba\u64fa\uff1b\u5982\u679c\u5176\u4e2d2\u4eba\u5404\u64fa4\u76c6\uff0c\u5176\u9918\u7684\u4eba\u5404\u64fa6\u76c6\uff0c\u9019\u4e9b\u82b1\u76c6\u6b63\u597d\u64fa\u5b8c\u3002\u554f\uff1a\u6709\u30107\u3011\u4f4d\u7ae5\u8ecd\u968a\u54e1\u3001\u5171\u64fa\u301038\u3011\u500b\u82b1\u76c6\", \"(21) \u5de5\u4eba\u92ea\u4e00\u689d\u8def\u57fa\uff0c\u82e5\u6bcf\u5929\u92ea260\u516c\u5c3a\uff0c\u92ea\u5b8c\u5168\u8def\u9577\u5c31\u5f97\u5ef6\u95778\u5929\uff1b\u82e5\u6bcf\u5929\u92ea300\u516c\u5c3a\uff0c\u92ea\u5b8c\u5168\u8def\u9577\u4ecd\u8981\u5ef6\u95774\u5929\u3002\u554f\uff1a\u9019\u689d\u8def\u9577\u30107800\u3011\u516c\u5c3a\u3002\", \"(22) \u4e00\u7fa4\u7334\u5b50\u5206\u6843\u5b50\uff0c\u5982\u679c\u6bcf\u96bb\u7334\u52065\u500b\uff0c\u9084\u991848\u500b\uff0c\u5982\u679c\u5176\u4e2d9\u96bb\u7334\u5404\u52066\u500b\u6843\uff0c\u5176\u9918\u7684\u7334\u52068\u500b\u6843\u5b50\uff0c\u6070\u597d\u5206\u5b8c\u3002\u554f\uff1a\u6709\u301022\u3011\u96bb\u7334\u5b50\u3001\u3010158\u3011\u6843\u5b50\", \"(23) \u5b78\u6821\u8cb7\u4f86\u4e00\u6279\u96fb\u98a8\u6247\u5206\u7d66\u5404\u73ed\u3002\u82e5\u6709\u5169\u500b\u73ed\u6bcf\u73ed\u5206\u52304\u53f0\uff0c\u5176\u9918\u6bcf\u73ed\u53ea\u80fd\u52062\u53f0\uff1b\u5982\u679c\u6709\u4e00\u500b\u73ed\u52066\u53f0\uff0c\u5176\u9918\u6bcf\u73ed\u52064\u53f0\uff0c\u9084\u5dee12\u53f0\u3002\u554f\uff1a\u5171\u8cb7\u4f86\u301018\u3011\u53f0\u96fb\u98a8\u6247\u3001\u6709\u30107\u3011\u500b\u73ed\", \"(24) \u5c0f\u570b\u8cb7\u4e86\u4e00\u672c\u300a\u8da3\u5473\u6578\u5b78\u300b\uff0c\u4ed6\u8a08\u756b\uff1a\u6bcf\u5929\u505a3\u984c\uff0c\u5247\u5269\u4e0b16\u984c\uff1b\u82e5\u6bcf\u5929\u505a5\u984c\uff0c\u5247\u6700\u5f8c\u4e00\u5929\u53ea\u8981\u505a1\u984c\u3002\u90a3\u9ebc\u9019\u672c\u66f8\u5171\u6709\u301046\u3011\u9053\u984c\uff1b\u5c0f\u570b\u8a08\u756b\u505a\u301010\u3011\u5929\", \"(25) \u4e09\u5e74\u7d1a\u7d66\u512a\u79c0\u5b78\u751f\u767c\u734e\u54c1\u66f8\uff0c\u5982\u679c\u6bcf\u500b\u5b78\u751f\u767c5\u518a\u9084\u526932\u518a\uff1b\u5982\u679c\u5176\u4e2d10\u500b\u5b78\u751f\u6bcf\u4eba\u767c4\u518a\uff0c\u5176\u9918\u6bcf\u4eba\u767c8\u518a\uff0c\u5c31\u6070\u597d\u767c\u5b8c\u3002\u554f\uff1a\u512a\u79c0\u5b78\u751f\u6709\u301024\u3011\u4eba\u3001\u734e\u54c1\u66f8\u6709\u3010152\u3011\u518a\", \"(26) \u6625\u7bc0\u524d\u5915\uff0c\u4e00\u500b\u5bcc\u7fc1\u65bd\u6368\u4e10\u5e6b\u5011\uff0c\u4e00\u958b\u59cb\u4ed6\u6e96\u5099\u7d66\u6bcf\u4eba100\u5143\uff0c\u7d50\u679c\u5269\u4e0b350\u5143\uff1b\u4ed6\u6c7a\u5b9a\u6bcf\u4eba\u591a\u7d6620\u5143\uff0c\u4f46\u9019\u6642\u5019\u8d95\u4f86\u4e865\u500b\u4e5e\u4e10\uff0c\u5982\u679c\u4ed6\u5011\u6bcf\u500b\u4eba\u62ff\u5230\u7684\u9322\u548c\u5176\u4ed6\u4e5e\u4e10\u4e00\u6a23\u591a\uff0c\u5bcc\u7fc1\u9084\u9700\u591a\u6e96\u5099550\u5143\u3002\u554f\uff1a\u539f\u4f86\u6709\u301015\u3011\u540d\u4e5e\u4e10\", \"(27) \u738b\u8001\u5e2b\u770b\u4e00\u672c\u5c0f\u8aaa\uff0c\u5982\u679c\u6bcf\u5929\u770b25\u9801\uff0c\u770b\u5b8c\u5168\u66f8\u6bd4\u898f\u5b9a\u6642\u9593\u591a\u4e00\u5929\uff0c\u5982\u679c\u6bcf\u5929\u770b30\u9801\u6700\u5f8c\u4e00\u5929\u5c11\u770b15\u9801\uff0c\u5982\u679c\u6bcf\u5929\u770b29\u9801\uff0c\u6700\u5f8c\u4e00\u5929\u8981\u8b80\u301022\u3011\u9801\u624d\u80fd\u6309\u898f\u5b9a\u8b80\u5b8c\", \"(28) \u8001\u5e2b\u628a\u4e00\u888b\u7cd6\u5206\u7d66\u5c0f\u670b\u53cb\u3002\u5982\u679c\u53ea\u5206\u7d66\u5c0f\u73ed\uff0c\u6bcf\u4eba\u53ef\u5f9712\u584a\uff0c\u5982\u679c\u5206\u7d66\u4e2d\u73ed\u548c\u5c0f\u73ed\uff0c\u6bcf\u4eba\u53ea\u80fd\u5206\u52304\u584a\u3002\u5982\u679c\u9019\u888b\u7cd6\u53ea\u5206\u7d66\u4e2d\u73ed\uff0c\u6bcf\u4eba\u53ef\u5206\u5230\u3010 6 \u3011\u584a\", \"(29) \u9ad4\u80b2\u968a\u5c07\u4e00\u4e9b\u7fbd\u6bdb\u7403\u5206\u7d66\u82e5\u5e72\u500b\u4eba\uff0c\u6bcf\u4eba5\u500b\u9084\u591a\u991810\u500b\u7fbd\u6bdb\u7403\uff0c\u5982\u679c\u4eba\u6578\u589e\u52a0\u5230 3\u500d\uff0c\u90a3\u9ebc\u6bcf\u4eba\u52062\u500b\u7fbd\u6bdb\u7403\u9084\u7f3a\u5c118\u500b\uff0c\u554f\uff1a\u6709\u7fbd\u6bdb\u7403\u3010100\u3011\u500b\", \"(30) \u516d\u5e74\u7d1a\u8209\u884c\u806f\u6b61\u665a\u6703\uff0c\u8001\u5e2b\u5e36\u8457\u4e00\u7b46\u9322\u53bb\u8cb7\u96f6\u98df\u3002\u5982\u679c\u8cb7\u7cd6\u679c13\u516c\u65a4\uff0c\u9084\u5dee4\u5143\uff1b\u5982\u679c\u8cb7\u725b\u5976\u7cd615\u516c\u65a4\uff0c\u5247\u9084\u52692\u5143\uff0c\u5df2\u77e5\u6bcf\u516c\u65a4\u7cd6\u679c\u6bd4\u725b\u5976\u7cd6\u8cb42\u5143\u3002\u554f\uff1a\u8001\u5e2b\u5e36\u3010 152 \u3011\u5143\", \"(31) 48\u672c\u66f8\u5206\u7d66\u5169\u7d44\u5c0f\u670b\u53cb\uff0c\u5df2\u77e5\u7b2c\u4e8c\u7d44\u6bd4\u7b2c\u4e00\u7d44\u591a5\u4eba\u3002\u5982\u679c\u628a\u66f8\u5168\u90e8\u5206\u7d66\u7b2c\u4e00\u7d44\uff0c\u90a3\u9ebc\u6bcf\u4eba4\u672c\uff0c\u6709\u5269\u9918\uff1b\u6bcf\u4eba5\u672c\uff0c\u66f8\u4e0d\u5920\u3002\u5982\u679c\u628a\u66f8\u5168\u5206\u7d66\u7b2c\u4e8c\u7d44\uff0c\u90a3\u9ebc\u6bcf\u4eba3\u672c\uff0c\u6709\u5269\u9918\uff1b\u6bcf\u4eba4\u672c\uff0c\u66f8\u4e0d\u5920\uff0c\u554f\uff1a\u7b2c\u4e00\u7d44\u6709\u301010\u3011\u4eba\u3001\u7b2c\u4e8c\u7d44\u6709\u301015\u3011\u4eba\", \"(32) \u4e00\u4e9b\u6854\u5b50\u5206\u7d66\u82e5\u5e72\u4eba\uff0c\u6bcf\u4eba5\u500b\u991810\u500b\u6854\u5b50\uff0c\u5982\u679c\u4eba\u6578\u589e\u52a0\u52303\u500d\u9084\u5c115\u4eba\uff0c\u90a3\u9ebc\u6bcf\u4eba\u52062\u500b\u9084\u7f3a8\u500b\uff0c\u6709\u6854\u5b50\u3010150\u3011\u500b\u3002\", \"(33) \u5e7c\u7a1a\u5712\u6559\u5e2b\u628a\u4e00\u7bb1\u9905\u4e7e\u5206\u7d66\u5c0f\u73ed\u548c\u4e2d\u73ed\u7684\u5c0f\u670b\u53cb\uff0c\u5e73\u5747\u6bcf\u4eba\u5206\u5f976\u584a\uff0c\u5982\u679c\u53ea\u5206\u7d66\u4e2d\u73ed\u5c0f\u670b\u53cb\uff0c\u5e73\u5747\u6bcf\u4eba\u53ef\u4ee5\u591a\u5206\u5f974\u584a\u3002\u554f\uff1a\u5982\u679c\u53ea\u5206\u7d66\u5c0f\u73ed\u7684\u5c0f\u670b\u53cb\uff0c\u5e73\u5747\u6bcf\u4eba\u5206\u5f97\u301015\u3011\u584a\", \"(34) \u8001\u5e2b\u628a\u4e00\u7c43\u860b\u679c\u5206\u7d66\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u6e1b\u5c11\u4e00\u540d\u540c\u5b78\uff0c\u6bcf\u500b\u540c\u5b78\u6b63\u597d\u5206\u5f975\u500b\uff1b\u5982\u679c\u589e\u52a0\u4e00\u540d\u540c\u5b78\uff0c\u6b63\u597d\u6bcf\u4eba\u5206\u5f974\u500b\u3002\u554f\uff1a\u5c0f\u670b\u53cb\u6709\u30109\u3011\u4eba\u3001\u860b\u679c\u6709\u301040\u3011\u500b\", \"(35) \u98df\u5802\u63a1\u8cfc\u54e1\u5c0f\u674e\u53bb\u8cb7\u8089\uff0c\u5982\u679c\u8cb7\u725b\u808918\u516c\u65a4\uff0c\u90a3\u9ebc\u5dee40\u5143\uff1b\u5982\u679c\u8cb7\u8c6c\u808920\u516c\u65a4\uff0c\u90a3\u9ebc\u591a20\u5143\u3002\u5df2\u77e5\u725b\u8089\u6bd4\u8c6c\u8089\u6bcf\u516c\u65a4\u8cb48\u5143\u3002\u554f\uff1a\u725b\u8089\u6bcf\u516c\u65a4\u301050\u3011\u5143\u3001\u8c6c\u8089\u6bcf\u516c\u65a4\u301042\u3011\u5143\uff0c\u5c0f\u674e\u5e36\u4e86\u3010860\u3011\u5143\", \"(36) \u56db\u5e74\u7d1a\u67d0\u73ed\u7684\u540c\u5b78\u53bb\u690d\u6a39\uff0c\u4ed6\u5011\u5206\u4e86\u4e00\u4e0b\u5c0f\u7d44\uff0c\u5982\u679c\u589e\u52a0\u4e00\u500b\u5c0f\u7d44\uff0c\u6b63\u597d\u6bcf\u5c0f\u7d445\u4eba\uff1b\u5982\u679c\u6e1b\u5c11\u4e00\u5c0f\u7d44\uff0c\u6b63\u597d\u6bcf\u7d447\u4eba\u3002\u554f\uff1a\u9019\u500b\u73ed\u5171\u6709\u301035\u3011\u4eba\", \"(37) \u7334\u738b\u5e36\u9818\u4e00\u7fa4\u7334\u5b50\u53bb\u6458\u6843\u3002\u4e0b\u5348\u6536\u5de5\u5f8c\uff0c\u7334\u738b\u958b\u59cb\u5206\u914d\uff0c\u82e5\u5927\u7334\u52065\u500b\uff0c\u5c0f\u7334\u52063\u500b\uff0c\u7334\u738b\u53ef\u755910\u500b\uff1b\u82e5\u5927\u3001\u5c0f\u7334\u90fd\u52064\u500b\uff0c\u7334\u738b\u80fd\u7559\u4e0b20\u500b\u3002\u5728\u9019\u7fa4\u7334\u5b50\u4e2d\uff0c\u5927\u7334\uff08\u4e0d\u5305\u62ec\u7334\u738b\uff09\u6bd4\u5c0f\u7334\u591a\u3010  10 \u3011\u96bb\u3002\", \"(38) \u5c0f\u660e\u5abd\u5abd\u5e36\u8457\u4e00\u7b46\u9322\u53bb\u8cb7\u8089\uff0c\u82e5\u8cb710\u516c\u65a4\u725b\u8089\u5247\u9084\u5dee6\u5143\uff0c\u82e5\u8cb712\u516c\u65a4\u8c6c\u8089\u5247\u9084\u52694\u5143\u3002\u5df2\u77e5\u6bcf\u516c\u65a4\u725b\u8089\u6bd4\u8c6c\u8089\u8cb43\u5143\uff0c\u554f\uff1a\u5c0f\u660e\u5abd\u5abd\u5e36\u4e86\u3010124\u3011\u5143\", \"(39) \u5e7c\u7a1a\u5712\u5c07\u4e00\u7b50\u860b\u679c\u5206\u7d66\u5927\u73ed\u548c\u5c0f\u73ed\u7684\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u5927\u73ed\u6bcf\u4eba\u52065\u500b\uff0c\u5c31\u591a10\u500b\uff1b\u5982\u679c\u5c0f\u73ed\u6bcf\u4eba\u52068\u500b\uff0c\u5c31\u5c11\u4e862\u500b\u3002\u5df2\u77e5\u5927\u73ed\u6bd4\u5c0f\u73ed\u591a3\u4eba\u3002\u554f\uff1a\u9019\u7b50\u860b\u679c\u6709\u301070\u3011\u500b\", \"(40) \u7532\u3001\u4e59\u5169\u7d44\u540c\u5b78\u505a\u7d05\u82b1\uff0c\u6bcf\u4eba\u505a8\u6735\uff0c\u6b63\u597d\u9001\u7d66\u4e94\u5e74\u7d1a\u6bcf\u500b\u540c\u5b78\u4e00\u6735\u3002\u5982\u679c\u628a\u9019\u4e9b\u7d05\u82b1\u8b93\u7532\u7d44\u55ae\u7368\u505a\uff0c\u6bcf\u4eba\u8981\u591a\u505a4\u6735\u3002\u5982\u679c\u628a\u9019\u4e9b\u7d05\u82b1\u8b93\u4e59\u7d44\u540c\u5b78\u55ae\u7368\u505a\uff0c\u6bcf\u4eba\u8981\u505a\u301024\u3011\u6735\"];\nconst TITLE = \"\u76c8\u8667\u554f\u984c - \u57fa\u672c\u578b\";\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\n// paras.items[0] = Title paragraph\n// paras.items[1..10] = existing question paragraphs (1)-(10)\nconst titlePara = paras.items[0];\ntitlePara.insertText(TITLE, Word.InsertLocation.replace);\n\n// Replace the text of the first 10 existing question paragraphs with the\n// first 10 new questions (keeps their pPr / jc=\"left\" formatting intact).\nfor (let i = 0; i < 10; i++) {\n  const p = paras.items[i + 1];\n  p.insertText(QUESTIONS[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// Append the remaining 30 questions (11)-(40) as new paragraphs after the\n// 10th question paragraph, chaining insertParagraph so each new paragraph\n// inherits the \"question\" style + left alignment from its predecessor.\nlet anchor = paras.items[10];\nfor (let i = 10; i < QUESTIONS.length; i++) {\n  anchor = anchor.insertParagraph(QUESTIONS[i], Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$titleText = \"\u76c8\u8667\u554f\u984c - \u57fa\u672c\u578b\"\n$d.Paragraphs.Item(1).Range.Text = $titleText\n\n$questions = @(\n  \"(1) \u738b\u8001\u5e2b\u53bb\u8cb7\u5152\u7ae5\u5c0f\u63d0\u7434\uff0c\u82e5\u8cb77\u628a\uff0c\u5247\u6240\u5e36\u7684\u9322\u5dee110\u5143\uff1b\u82e5\u8cb75\u628a\uff0c\u5247\u6240\u5e36\u7684\u9322\u9084\u5dee30\u5143\u3002\u554f\uff1a\u5152\u7ae5\u5c0f\u63d0\u7434\u4e00\u628a\u301040\u3011\u5143\uff0c\u738b\u8001\u5e2b\u5e36\u4e86\u3010170\u3011\u5143\u3002\",\n  \"(2) \u4e00\u500b\u5546\u8ca9\u4f30\u8a08\uff0c\u82e51\u516c\u65a4\u860b\u679c\u8ce324\u5143\uff0c\u6703\u8ce040\u5143\uff1b\u82e51\u516c\u65a4\u860b\u679c\u8ce330\u5143\uff0c\u53ef\u4ee5\u8cfa80\u5143\u3002\u554f\uff1a\u82e5\u4ee5\u4e0d\u8ce0\u4e0d\u8cfa\u7684\u50f9\u683c\u8ce3\u51fa\uff0c\u6bcf\u516c\u65a4\u860b\u679c\u61c9\u8ce3\u301026\u3011\u5143\",\n  \"(3) \u8c93\u5abd\u5abd\u7d66\u5c0f\u8c93\u5206\u9b5a\uff0c\u6bcf\u96bb\u5c0f\u8c93\u520610\u689d\u9b5a\uff0c\u5c31\u591a\u51fa8\u689d\u9b5a\uff0c\u6bcf\u96bb\u5c0f\u8c93\u520611\u689d\u9b5a\u5247\u6b63\u597d\u5206\u5b8c\u3002\u554f\uff1a\u4e00\u5171\u6709\u30108\u3011\u96bb\u5c0f\u8c93\u3001\u8c93\u5abd\u5abd\u4e00\u5171\u6709\u301088\u3011\u689d\u9b5a\",\n  \"(4) \u5b78\u751f\u5011\u53c3\u52a0\u690d\u6a39\u6d3b\u52d5\uff0c\u5982\u679c\u6bcf\u4eba\u683d5\u68f5\u6a39\uff0c\u9084\u526912\u68f5\u6a39\uff1b\u5982\u679c\u6bcf\u4eba\u683d7\u68f5\uff0c\u5c31\u7f3a4\u68f5\u3002\u554f\uff1a\u5b78\u751f\u6709\u30108\u3011\u4eba\u3001\u4e00\u5171\u8981\u683d\u301052\u3011\u68f5\u6a39\",\n  \"(5) 401\u73ed\u540c\u5b78\u690d\u6a39\uff0c\u6bcf\u4eba\u690d1\u68f5\u9084\u526920\u68f5\uff0c\u6bcf\u4eba\u690d2\u68f5\u5dee30\u68f5\u3002\u554f\uff1a\u6709\u301050\u3011\u500b\u540c\u5b78\u3001\u301070\u3011\u68f5\u6a39\u82d7\",\n  \"(6) \u5abd\u5abd\u8cb7\u56de\u4e00\u7b50\u860b\u679c\uff0c\u6309\u8a08\u5283\u5403\u7684\u5929\u6578\u7b97\u4e86\u4e00\u4e0b\uff0c\u5982\u679c\u6bcf\u5929\u54034\u500b\uff0c\u8981\u591a\u51fa48\u500b\u860b\u679c\uff1b\u5982\u679c\u6bcf\u5929\u54036\u500b\uff0c\u5247\u53c8\u5c118\u500b\u860b\u679c\u3002\u554f\uff1a\u5abd\u5abd\u8cb7\u56de\u7684\u860b\u679c\u6709\u3010160\u3011\u500b\u3001\u8a08\u756b\u5403\u301028\u3011\u5929 \",\n  \"(7) \u5c0f\u8ecd\u5c07\u81ea\u5df1\u6536\u85cf\u7684\u4e00\u4e9b\u7167\u7247\u9001\u7d66\u5e7c\u5152\u5712\u5927\u73ed\u7684\u5c0f\u670b\u53cb\u5011\u3002\u5982\u679c\u6bcf\u4eba\u52069\u5f35\uff0c\u9084\u591a12\u5f35\uff0c\u5982\u679c\u6bcf\u4eba\u520610\u5f35\u5247\u6b63\u597d\u5206\u5b8c\u3002\u554f\uff1a\u6709\u301012\u3011\u500b\u5c0f\u670b\u53cb\uff0c\u7167\u7247\u6709\u3010120\u3011\u5f35\",\n  \"(8) \u4e00\u5305\u7cd6\u5206\u7d66\u5e7e\u500b\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u6bcf\u4eba\u52063\u584a\uff0c\u5247\u99183\u584a\uff1b\u5982\u679c\u6bcf\u4eba\u52065\u584a\uff0c\u5247\u5c111\u584a\u3002\u554f\uff1a\u5c0f\u670b\u53cb\u6709\u30102\u3011\u4eba\u3001\u7cd6\u6709\u30109\u3011\u584a\",\n  \"(9) \u7528\u7b50\u88dd\u897f\u74dc\uff0c\u5982\u679c\u6bcf\u7b50\u88dd5\u500b\uff0c\u5247\u5c1115\u500b\u897f\u74dc\uff1b\u5982\u679c\u6bcf\u7b50\u88dd3\u500b\uff0c\u5247\u591a29\u500b\u897f\u74dc\u3002\u554f\uff1a\u5171\u6709\u7b50\u301022\u3011\u500b\u3001\u897f\u74dc\u301095\u3011\u500b\",\n  \"(10) \u6709\u4e00\u6279\u7df4\u7fd2\u672c\u767c\u7d66\u5b78\u751f\uff0c\u5982\u679c\u6bcf\u4eba5\u672c\uff0c\u5247\u591a70\u672c\uff0c\u5982\u679c\u6bcf\u4eba7\u672c\uff0c\u5247\u591a10\u672c\u3002\u554f\uff1a\u9019\u500b\u73ed\u6709\u5b78\u751f\u301030\u3011\u4eba\uff0c\u6709\u7df4\u7fd2\u672c\u3010220\u3011\u672c\u3002\",\n  \"(11) \u8001\u7334\u5b50\u7d66\u5c0f\u7334\u5b50\u5206\u6843\uff0c\u6bcf\u96bb\u5c0f\u7334\u520610\u500b\u6843\uff0c\u5c31\u591a\u51fa9\u500b\u6843\uff0c\u6bcf\u96bb\u5c0f\u7334\u520611\u500b\u6843\u5247\u591a\u51fa2\u500b\u6843\u3002\u554f\uff1a\u4e00\u5171\u6709\u3010 7 \u3011\u96bb\u5c0f\u7334\u5b50\uff0c\u8001\u7334\u5b50\u4e00\u5171\u6709\u3010 79 \u3011\u500b\u6843\u5b50\",\n  \"(12) \u5b78\u6821\u5c07\u4e00\u6279\u925b\u7b46\u734e\u7d66\u5b78\u751f\uff0c\u6bcf\u4eba9\u652f\u7f3a15\u652f\uff1b\u6bcf\u4eba7\u652f\u7f3a7\u652f\u3002\u554f\uff1a\u5b78\u751f\u6709\u30104\u3011\u4eba\uff0c\u925b\u7b46\u6709\u301021\u3011\u652f\",\n  \"(13) \u5c0f\u670b\u53cb\u5206\u7cd6\u679c\uff0c\u6bcf\u4eba3\u7c92\uff0c\u991830\u7c92\uff1b\u6bcf\u4eba5\u7c92\uff0c\u5c114\u7c92\u3002\u554f\uff1a\u6709\u301017\u3011\u500b\u5c0f\u670b\u53cb\u3001\u301081\u3011\u7c92\u7cd6\",\n  \"(14) \u4e00\u500b\u6c7d\u8eca\u968a\u904b\u8f38\u4e00\u6279\u8ca8\u7269\uff0c\u5982\u679c\u6bcf\u8f1b\u6c7d\u8eca\u904b3500\u516c\u65a4\uff0c\u90a3\u9ebc\u8ca8\u7269\u9084\u5269\u4e0b5000\u516c\u65a4\uff1b\u5982\u679c\u6bcf\u8f1b\u6c7d\u8eca\u904b4000\u516c\u65a4\uff0c\u90a3\u9ebc\u8ca8\u7269\u9084\u5269\u4e0b500\u516c\u65a4\u3002\u554f\uff1a\u9019\u500b\u6c7d\u8eca\u968a\u6709\u30109\u3011\u8f1b\u6c7d\u8eca\u3001\u8981\u904b\u7684\u8ca8\u7269\u6709\u301036500\u3011\u516c\u65a4\",\n  \"(15) \u5b78\u6821\u5206\u914d\u82e5\u5e72\u4eba\u64e6\u73bb\u7483\uff0c\u5176\u4e2d\u5169\u4eba\u5404\u64e64\u584a\uff0c\u5176\u9918\u5404\u64e65\u584a\uff0c\u5247\u991812\u584a\uff1b\u82e5\u6bcf\u4eba\u64e66\u584a\uff0c\u5247\u6b63\u597d\u64e6\u5b8c\u3002\u554f\uff1a\u64e6\u73bb\u7483\u7684\u6709\u301010\u3011\u4eba\u3001\u73bb\u7483\u6709\u301060\u3011\u584a\",\n  \"(16) \u5c07\u4e00\u5806\u6843\u5b50\uff0c\u5e73\u5747\u5206\u7d66\u7334\u5b50\uff0c\u6bcf\u96bb\u7334\u5b50\u520610\u500b\uff0c\u6709\u5169\u96bb\u7334\u5b50\u6c92\u5206\u5230\uff0c\u7b2c\u4e8c\u6b21\u91cd\u5206\uff0c\u6bcf\u96bb\u7334\u5b508\u500b\u6843\u5b50\uff0c\u525b\u5de7\u5206\u5b8c\u3002\u554f\uff1a\u4e00\u5806\u6843\u5b50\u6709\u301080\u3011\u500b\u3001\u7334\u5b50\u6709\u301010\u3011\u96bb\",\n  \"(17) \u8fb2\u6c11\u92e4\u8349\uff0c\u5176\u4e2d5\u4eba\u5404\u92e44\u755d\uff0c\u9918\u4e0b\u7684\u5404\u92e43\u755d\uff0c\u9019\u6a23\u5206\u914d\u6700\u5f8c\u9918\u4e0b26\u755d\uff1b\u5982\u679c\u5176\u4e2d3\u4eba\u6bcf\u4eba\u5404\u92e43\u755d\uff0c\u9918\u4e0b\u7684\u4eba\u5404\u92e45\u755d\uff0c\u6700\u5f8c\u9918\u4e0b3\u755d\u3002\u554f\uff1a\u8349\u5730\u9762\u7a4d\u301082\u3011\u755d\u3001\u92e4\u8349\u4eba\u6578\u301017\u3011\u4eba\",\n  \"(18) \u8001\u5e2b\u628a\u4e00\u7c43\u860b\u679c\u5206\u7d66\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u6e1b\u5c11\u4e00\u540d\u540c\u5b78\uff0c\u6bcf\u500b\u540c\u5b78\u6b63\u597d\u5206\u5f975\u500b\uff1b\u5982\u679c\u589e\u52a0\u4e00\u500b\u540c\u5b78\uff0c\u6b63\u597d\u6bcf\u4eba\u5206\u5f974\u500b\uff0c\u6c42\u9019\u7c43\u860b\u679c\u4e00\u5171\u6709\u3010 40 \u3011\u500b\",\n  \"(19) \u5321\u660e\u548c\u674e\u6b23\u7d04\u5b9a\u5728\u76f8\u540c\u7684\u6642\u9593\u5167\u505a\u5b8c\u540c\u6a23\u7684\u4e00\u672c\u6578\u5b78\u7df4\u7fd2\u984c\u3002\u5321\u660e\u8a08\u756b\u982d\u5169\u5468\u6bcf\u9031\u505a30\u9053\uff0c\u4ee5\u5f8c\u6bcf\u9031\u505a25\u9053\uff1b\u674e\u6b23\u8a08\u756b\u982d\u5169\u5468\u6bcf\u9031\u505a35\u9053\uff0c\u4ee5\u5f8c\u6bcf\u9031\u505a30\u9053\u3002\u7d50\u679c\u674e\u6b23\u63d0\u524d\u5169\u5468\u505a\u5b8c\u3002\u554f\uff1a\u4ed6\u5011\u6e96\u5099\u301012\u3011\u5468\u505a\u5b8c\uff0c\u9019\u672c\u7df4\u7fd2\u984c\u5171\u6709\u3010310\u3011\u984c\",\n  \"(20) \u7ae5\u8ecd\u968a\u54e1\u64fa\u82b1\u76c6\u5e03\u7f6e\u6821\u5712\u3002\u5982\u679c\u6bcf\u4eba\u64fa5\u76c6\u82b1\uff0c\u9084\u67093\u76c6\u6c92\u4eba\u64fa\uff1b\u5982\u679c\u5176\u4e2d2\u4eba\u5404\u64fa4\u76c6\uff0c\u5176\u9918\u7684\u4eba\u5404\u64fa6\u76c6\uff0c\u9019\u4e9b\u82b1\u76c6\u6b63\u597d\u64fa\u5b8c\u3002\u554f\uff1a\u6709\u30107\u3011\u4f4d\u7ae5\u8ecd\u968a\u54e1\u3001\u5171\u64fa\u301038\u3011\u500b\u82b1\u76c6\",\n  \"(21) \u5de5\u4eba\u92ea\u4e00\u689d\u8def\u57fa\uff0c\u82e5\u6bcf\u5929\u92ea260\u516c\u5c3a\uff0c\u92ea\u5b8c\u5168\u8def\u9577\u5c31\u5f97\u5ef6\u95778\u5929\uff1b\u82e5\u6bcf\u5929\u92ea300\u516c\u5c3a\uff0c\u92ea\u5b8c\u5168\u8def\u9577\u4ecd\u8981\u5ef6\u95774\u5929\u3002\u554f\uff1a\u9019\u689d\u8def\u9577\u30107800\u3011\u516c\u5c3a\u3002\",\n  \"(22) \u4e00\u7fa4\u7334\u5b50\u5206\u6843\u5b50\uff0c\u5982\u679c\u6bcf\u96bb\u7334\u52065\u500b\uff0c\u9084\u991848\u500b\uff0c\u5982\u679c\u5176\u4e2d9\u96bb\u7334\u5404\u52066\u500b\u6843\uff0c\u5176\u9918\u7684\u7334\u52068\u500b\u6843\u5b50\uff0c\u6070\u597d\u5206\u5b8c\u3002\u554f\uff1a\u6709\u301022\u3011\u96bb\u7334\u5b50\u3001\u3010158\u3011\u6843\u5b50\",\n  \"(23) \u5b78\u6821\u8cb7\u4f86\u4e00\u6279\u96fb\u98a8\u6247\u5206\u7d66\u5404\u73ed\u3002\u82e5\u6709\u5169\u500b\u73ed\u6bcf\u73ed\u5206\u52304\u53f0\uff0c\u5176\u9918\u6bcf\u73ed\u53ea\u80fd\u52062\u53f0\uff1b\u5982\u679c\u6709\u4e00\u500b\u73ed\u52066\u53f0\uff0c\u5176\u9918\u6bcf\u73ed\u52064\u53f0\uff0c\u9084\u5dee12\u53f0\u3002\u554f\uff1a\u5171\u8cb7\u4f86\u301018\u3011\u53f0\u96fb\u98a8\u6247\u3001\u6709\u30107\u3011\u500b\u73ed\",\n  \"(24) \u5c0f\u570b\u8cb7\u4e86\u4e00\u672c\u300a\u8da3\u5473\u6578\u5b78\u300b\uff0c\u4ed6\u8a08\u756b\uff1a\u6bcf\u5929\u505a3\u984c\uff0c\u5247\u5269\u4e0b16\u984c\uff1b\u82e5\u6bcf\u5929\u505a5\u984c\uff0c\u5247\u6700\u5f8c\u4e00\u5929\u53ea\u8981\u505a1\u984c\u3002\u90a3\u9ebc\u9019\u672c\u66f8\u5171\u6709\u301046\u3011\u9053\u984c\uff1b\u5c0f\u570b\u8a08\u756b\u505a\u301010\u3011\u5929\",\n  \"(25) \u4e09\u5e74\u7d1a\u7d66\u512a\u79c0\u5b78\u751f\u767c\u734e\u54c1\u66f8\uff0c\u5982\u679c\u6bcf\u500b\u5b78\u751f\u767c5\u518a\u9084\u526932\u518a\uff1b\u5982\u679c\u5176\u4e2d10\u500b\u5b78\u751f\u6bcf\u4eba\u767c4\u518a\uff0c\u5176\u9918\u6bcf\u4eba\u767c8\u518a\uff0c\u5c31\u6070\u597d\u767c\u5b8c\u3002\u554f\uff1a\u512a\u79c0\u5b78\u751f\u6709\u301024\u3011\u4eba\u3001\u734e\u54c1\u66f8\u6709\u3010152\u3011\u518a\",\n  \"(26) \u6625\u7bc0\u524d\u5915\uff0c\u4e00\u500b\u5bcc\u7fc1\u65bd\u6368\u4e10\u5e6b\u5011\uff0c\u4e00\u958b\u59cb\u4ed6\u6e96\u5099\u7d66\u6bcf\u4eba100\u5143\uff0c\u7d50\u679c\u5269\u4e0b350\u5143\uff1b\u4ed6\u6c7a\u5b9a\u6bcf\u4eba\u591a\u7d6620\u5143\uff0c\u4f46\u9019\u6642\u5019\u8d95\u4f86\u4e865\u500b\u4e5e\u4e10\uff0c\u5982\u679c\u4ed6\u5011\u6bcf\u500b\u4eba\u62ff\u5230\u7684\u9322\u548c\u5176\u4ed6\u4e5e\u4e10\u4e00\u6a23\u591a\uff0c\u5bcc\u7fc1\u9084\u9700\u591a\u6e96\u5099550\u5143\u3002\u554f\uff1a\u539f\u4f86\u6709\u301015\u3011\u540d\u4e5e\u4e10\",\n  \"(27) \u738b\u8001\u5e2b\u770b\u4e00\u672c\u5c0f\u8aaa\uff0c\u5982\u679c\u6bcf\u5929\u770b25\u9801\uff0c\u770b\u5b8c\u5168\u66f8\u6bd4\u898f\u5b9a\u6642\u9593\u591a\u4e00\u5929\uff0c\u5982\u679c\u6bcf\u5929\u770b30\u9801\u6700\u5f8c\u4e00\u5929\u5c11\u770b15\u9801\uff0c\u5982\u679c\u6bcf\u5929\u770b29\u9801\uff0c\u6700\u5f8c\u4e00\u5929\u8981\u8b80\u301022\u3011\u9801\u624d\u80fd\u6309\u898f\u5b9a\u8b80\u5b8c\",\n  \"(28) \u8001\u5e2b\u628a\u4e00\u888b\u7cd6\u5206\u7d66\u5c0f\u670b\u53cb\u3002\u5982\u679c\u53ea\u5206\u7d66\u5c0f\u73ed\uff0c\u6bcf\u4eba\u53ef\u5f9712\u584a\uff0c\u5982\u679c\u5206\u7d66\u4e2d\u73ed\u548c\u5c0f\u73ed\uff0c\u6bcf\u4eba\u53ea\u80fd\u5206\u52304\u584a\u3002\u5982\u679c\u9019\u888b\u7cd6\u53ea\u5206\u7d66\u4e2d\u73ed\uff0c\u6bcf\u4eba\u53ef\u5206\u5230\u3010 6 \u3011\u584a\",\n  \"(29) \u9ad4\u80b2\u968a\u5c07\u4e00\u4e9b\u7fbd\u6bdb\u7403\u5206\u7d66\u82e5\u5e72\u500b\u4eba\uff0c\u6bcf\u4eba5\u500b\u9084\u591a\u991810\u500b\u7fbd\u6bdb\u7403\uff0c\u5982\u679c\u4eba\u6578\u589e\u52a0\u5230 3\u500d\uff0c\u90a3\u9ebc\u6bcf\u4eba\u52062\u500b\u7fbd\u6bdb\u7403\u9084\u7f3a\u5c118\u500b\uff0c\u554f\uff1a\u6709\u7fbd\u6bdb\u7403\u3010100\u3011\u500b\",\n  \"(30) \u516d\u5e74\u7d1a\u8209\u884c\u806f\u6b61\u665a\u6703\uff0c\u8001\u5e2b\u5e36\u8457\u4e00\u7b46\u9322\u53bb\u8cb7\u96f6\u98df\u3002\u5982\u679c\u8cb7\u7cd6\u679c13\u516c\u65a4\uff0c\u9084\u5dee4\u5143\uff1b\u5982\u679c\u8cb7\u725b\u5976\u7cd615\u516c\u65a4\uff0c\u5247\u9084\u52692\u5143\uff0c\u5df2\u77e5\u6bcf\u516c\u65a4\u7cd6\u679c\u6bd4\u725b\u5976\u7cd6\u8cb42\u5143\u3002\u554f\uff1a\u8001\u5e2b\u5e36\u3010 152 \u3011\u5143\",\n  \"(31) 48\u672c\u66f8\u5206\u7d66\u5169\u7d44\u5c0f\u670b\u53cb\uff0c\u5df2\u77e5\u7b2c\u4e8c\u7d44\u6bd4\u7b2c\u4e00\u7d44\u591a5\u4eba\u3002\u5982\u679c\u628a\u66f8\u5168\u90e8\u5206\u7d66\u7b2c\u4e00\u7d44\uff0c\u90a3\u9ebc\u6bcf\u4eba4\u672c\uff0c\u6709\u5269\u9918\uff1b\u6bcf\u4eba5\u672c\uff0c\u66f8\u4e0d\u5920\u3002\u5982\u679c\u628a\u66f8\u5168\u5206\u7d66\u7b2c\u4e8c\u7d44\uff0c\u90a3\u9ebc\u6bcf\u4eba3\u672c\uff0c\u6709\u5269\u9918\uff1b\u6bcf\u4eba4\u672c\uff0c\u66f8\u4e0d\u5920\uff0c\u554f\uff1a\u7b2c\u4e00\u7d44\u6709\u301010\u3011\u4eba\u3001\u7b2c\u4e8c\u7d44\u6709\u301015\u3011\u4eba\",\n  \"(32) \u4e00\u4e9b\u6854\u5b50\u5206\u7d66\u82e5\u5e72\u4eba\uff0c\u6bcf\u4eba5\u500b\u991810\u500b\u6854\u5b50\uff0c\u5982\u679c\u4eba\u6578\u589e\u52a0\u52303\u500d\u9084\u5c115\u4eba\uff0c\u90a3\u9ebc\u6bcf\u4eba\u52062\u500b\u9084\u7f3a8\u500b\uff0c\u6709\u6854\u5b50\u3010150\u3011\u500b\u3002\",\n  \"(33) \u5e7c\u7a1a\u5712\u6559\u5e2b\u628a\u4e00\u7bb1\u9905\u4e7e\u5206\u7d66\u5c0f\u73ed\u548c\u4e2d\u73ed\u7684\u5c0f\u670b\u53cb\uff0c\u5e73\u5747\u6bcf\u4eba\u5206\u5f976\u584a\uff0c\u5982\u679c\u53ea\u5206\u7d66\u4e2d\u73ed\u5c0f\u670b\u53cb\uff0c\u5e73\u5747\u6bcf\u4eba\u53ef\u4ee5\u591a\u5206\u5f974\u584a\u3002\u554f\uff1a\u5982\u679c\u53ea\u5206\u7d66\u5c0f\u73ed\u7684\u5c0f\u670b\u53cb\uff0c\u5e73\u5747\u6bcf\u4eba\u5206\u5f97\u301015\u3011\u584a\",\n  \"(34) \u8001\u5e2b\u628a\u4e00\u7c43\u860b\u679c\u5206\u7d66\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u6e1b\u5c11\u4e00\u540d\u540c\u5b78\uff0c\u6bcf\u500b\u540c\u5b78\u6b63\u597d\u5206\u5f975\u500b\uff1b\u5982\u679c\u589e\u52a0\u4e00\u540d\u540c\u5b78\uff0c\u6b63\u597d\u6bcf\u4eba\u5206\u5f974\u500b\u3002\u554f\uff1a\u5c0f\u670b\u53cb\u6709\u30109\u3011\u4eba\u3001\u860b\u679c\u6709\u301040\u3011\u500b\",\n  \"(35) \u98df\u5802\u63a1\u8cfc\u54e1\u5c0f\u674e\u53bb\u8cb7\u8089\uff0c\u5982\u679c\u8cb7\u725b\u808918\u516c\u65a4\uff0c\u90a3\u9ebc\u5dee40\u5143\uff1b\u5982\u679c\u8cb7\u8c6c\u808920\u516c\u65a4\uff0c\u90a3\u9ebc\u591a20\u5143\u3002\u5df2\u77e5\u725b\u8089\u6bd4\u8c6c\u8089\u6bcf\u516c\u65a4\u8cb48\u5143\u3002\u554f\uff1a\u725b\u8089\u6bcf\u516c\u65a4\u301050\u3011\u5143\u3001\u8c6c\u8089\u6bcf\u516c\u65a4\u301042\u3011\u5143\uff0c\u5c0f\u674e\u5e36\u4e86\u3010860\u3011\u5143\",\n  \"(36) \u56db\u5e74\u7d1a\u67d0\u73ed\u7684\u540c\u5b78\u53bb\u690d\u6a39\uff0c\u4ed6\u5011\u5206\u4e86\u4e00\u4e0b\u5c0f\u7d44\uff0c\u5982\u679c\u589e\u52a0\u4e00\u500b\u5c0f\u7d44\uff0c\u6b63\u597d\u6bcf\u5c0f\u7d445\u4eba\uff1b\u5982\u679c\u6e1b\u5c11\u4e00\u5c0f\u7d44\uff0c\u6b63\u597d\u6bcf\u7d447\u4eba\u3002\u554f\uff1a\u9019\u500b\u73ed\u5171\u6709\u301035\u3011\u4eba\",\n  \"(37) \u7334\u738b\u5e36\u9818\u4e00\u7fa4\u7334\u5b50\u53bb\u6458\u6843\u3002\u4e0b\u5348\u6536\u5de5\u5f8c\uff0c\u7334\u738b\u958b\u59cb\u5206\u914d\uff0c\u82e5\u5927\u7334\u52065\u500b\uff0c\u5c0f\u7334\u52063\u500b\uff0c\u7334\u738b\u53ef\u755910\u500b\uff1b\u82e5\u5927\u3001\u5c0f\u7334\u90fd\u52064\u500b\uff0c\u7334\u738b\u80fd\u7559\u4e0b20\u500b\u3002\u5728\u9019\u7fa4\u7334\u5b50\u4e2d\uff0c\u5927\u7334\uff08\u4e0d\u5305\u62ec\u7334\u738b\uff09\u6bd4\u5c0f\u7334\u591a\u3010  10 \u3011\u96bb\u3002\",\n  \"(38) \u5c0f\u660e\u5abd\u5abd\u5e36\u8457\u4e00\u7b46\u9322\u53bb\u8cb7\u8089\uff0c\u82e5\u8cb710\u516c\u65a4\u725b\u8089\u5247\u9084\u5dee6\u5143\uff0c\u82e5\u8cb712\u516c\u65a4\u8c6c\u8089\u5247\u9084\u52694\u5143\u3002\u5df2\u77e5\u6bcf\u516c\u65a4\u725b\u8089\u6bd4\u8c6c\u8089\u8cb43\u5143\uff0c\u554f\uff1a\u5c0f\u660e\u5abd\u5abd\u5e36\u4e86\u3010124\u3011\u5143\",\n  \"(39) \u5e7c\u7a1a\u5712\u5c07\u4e00\u7b50\u860b\u679c\u5206\u7d66\u5927\u73ed\u548c\u5c0f\u73ed\u7684\u5c0f\u670b\u53cb\uff0c\u5982\u679c\u5927\u73ed\u6bcf\u4eba\u52065\u500b\uff0c\u5c31\u591a10\u500b\uff1b\u5982\u679c\u5c0f\u73ed\u6bcf\u4eba\u52068\u500b\uff0c\u5c31\u5c11\u4e862\u500b\u3002\u5df2\u77e5\u5927\u73ed\u6bd4\u5c0f\u73ed\u591a3\u4eba\u3002\u554f\uff1a\u9019\u7b50\u860b\u679c\u6709\u301070\u3011\u500b\",\n  \"(40) \u7532\u3001\u4e59\u5169\u7d44\u540c\u5b78\u505a\u7d05\u82b1\uff0c\u6bcf\u4eba\u505a8\u6735\uff0c\u6b63\u597d\u9001\u7d66\u4e94\u5e74\u7d1a\u6bcf\u500b\u540c\u5b78\u4e00\u6735\u3002\u5982\u679c\u628a\u9019\u4e9b\u7d05\u82b1\u8b93\u7532\u7d44\u55ae\u7368\u505a\uff0c\u6bcf\u4eba\u8981\u591a\u505a4\u6735\u3002\u5982\u679c\u628a\u9019\u4e9b\u7d05\u82b1\u8b93\u4e59\u7d44\u540c\u5b78\u55ae\u7368\u505a\uff0c\u6bcf\u4eba\u8981\u505a\u301024\u3011\u6735\"\n)\n\n# Replace the text of the first 10 existing question paragraphs (index 2..11)\nfor ($i = 0; $i -lt 10; $i++) {\n  $d.Paragraphs.Item($i + 2).Range.Text = $questions[$i]\n}\n\n# Append the remaining 30 questions (11)-(40) as new paragraphs, chaining\n# InsertParagraphAfter so each new paragraph inherits the \"question\" style\n# and left alignment from its predecessor.\n$anchor = $d.Paragraphs.Item(11)\nfor ($i = 10; $i -lt $questions.Count; $i++) {\n  $anchor.Range.InsertParagraphAfter()\n  $newIndex = $anchor.Index + 1\n  $newPara = $d.Paragraphs.Item($newIndex)\n  $newPara.Range.Text = $questions[$i]\n  $anchor = $newPara\n}\n"}
